$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.17606616307731
$ws.Range("C2").Value = 7.437949852043538
$ws.Range("D2").Value = 2.817955889696813
$ws.Range("E2").Value = 4.68286913085535
$ws.Range("F2").Value = 45.28983511499318
$ws.Range("G2").Value = 2.169427264939894
$ws.Range("H2").Value = 4.918938535312271
$ws.Range("I2").Value = 5.552781683948353
$ws.Range("K2").Value = 31.83522731462135
$ws.Range("L2").Value = 5.574403303264908
$ws.Range("M2").Value = 11.68515567508415
$ws.Range("N2").Value = 6.37663000761913
$ws.Range("B3").Value = 12.42160090677149
$ws.Range("C3").Value = 7.000246700013244
$ws.Range("D3").Value = 2.830154488745834
$ws.Range("E3").Value = 4.518593133140382
$ws.Range("F3").Value = 43.80874731094075
$ws.Range("G3").Value = 2.175217529732933
$ws.Range("H3").Value = 5.202794239097274
$ws.Range("I3").Value = 5.789279069321565
$ws.Range("K3").Value = 30.96407826341639
$ws.Range("L3").Value = 5.464644338715519
$ws.Range("M3").Value = 11.10211760429731
$ws.Range("N3").Value = 6.182204476008822
$ws.Range("B4").Value = 11.93709125354466
$ws.Range("C4").Value = 6.726843973495944
$ws.Range("D4").Value = 2.836423294488264
$ws.Range("E4").Value = 4.414734246523651
$ws.Range("F4").Value = 42.86907852103264
$ws.Range("G4").Value = 2.178886047003021
$ws.Range("H4").Value = 5.382947573749769
$ws.Range("I4").Value = 5.939775645621434
$ws.Range("K4").Value = 30.41073478824639
$ws.Range("L4").Value = 5.394742859896044
$ws.Range("M4").Value = 10.73298411230097
$ws.Range("N4").Value = 6.060890055636256
$ws.Range("B5").Value = 11.73358397154166
$ws.Range("C5").Value = 6.620213421354597
$ws.Range("D5").Value = 2.837217295429793
$ws.Range("E5").Value = 4.371847782255514
$ws.Range("F5").Value = 42.44677567198161
$ws.Range("G5").Value = 2.180423561641159
$ws.Range("H5").Value = 5.458531033824801
$ws.Range("I5").Value = 6.004805702151616
$ws.Range("K5").Value = 30.15634631627625
$ws.Range("L5").Value = 5.364205341814371
$ws.Range("M5").Value = 10.58127528054277
$ws.Range("N5").Value = 6.01274345810149
$ws.Range("B6").Value = 11.69842104761963
$ws.Range("C6").Value = 6.609561740489406
$ws.Range("D6").Value = 2.835572329261034
$ws.Range("E6").Value = 4.364908493440372
$ws.Range("F6").Value = 42.33688462702521
$ws.Range("G6").Value = 2.180697041474108
$ws.Range("H6").Value = 5.471984745253407
$ws.Range("I6").Value = 6.018646729108176
$ws.Range("K6").Value = 30.08392664653353
$ws.Range("L6").Value = 5.357360186807531
$ws.Range("M6").Value = 10.55756583648104
$ws.Range("N6").Value = 6.006843804708244
$ws.Range("B7").Value = 11.9314900825083
$ws.Range("C7").Value = 6.744409652811104
$ws.Range("D7").Value = 2.831642402907002
$ws.Range("E7").Value = 4.414769702693364
$ws.Range("F7").Value = 42.7571869860182
$ws.Range("G7").Value = 2.178951019006887
$ws.Range("H7").Value = 5.386156118661327
$ws.Range("I7").Value = 5.948522422438883
$ws.Range("K7").Value = 30.32638177332068
$ws.Range("L7").Value = 5.389631770621195
$ws.Range("M7").Value = 10.73529939519641
$ws.Range("N7").Value = 6.065972355405987
$ws.Range("B8").Value = 12.91671504710697
$ws.Range("C8").Value = 7.311366517617008
$ws.Range("D8").Value = 2.816244277135253
$ws.Range("E8").Value = 4.627617864544329
$ws.Range("F8").Value = 44.65082214843456
$ws.Range("G8").Value = 2.171457430801771
$ws.Range("H8").Value = 5.018427156232992
$ws.Range("I8").Value = 5.643235566621291
$ws.Range("K8").Value = 31.43536396345024
$ws.Range("L8").Value = 5.531042705320422
$ws.Range("M8").Value = 11.49181017875025
$ws.Range("N8").Value = 6.317329817957503
$ws.Range("B9").Value = 14.68006849616475
$ws.Range("C9").Value = 8.474007024264541
$ws.Range("D9").Value = 2.784269415288499
$ws.Range("E9").Value = 5.020199770503932
$ws.Range("F9").Value = 48.28769271949385
$ws.Range("G9").Value = 2.15755498597489
$ws.Range("H9").Value = 4.33867414326725
$ws.Range("I9").Value = 5.07235578876985
$ws.Range("K9").Value = 33.59316987765453
$ws.Range("L9").Value = 5.79655743344706
$ws.Range("M9").Value = 12.87354760836024
$ws.Range("N9").Value = 6.783344631682027
$ws.Range("B10").Value = 15.86530003514611
$ws.Range("C10").Value = 9.310349645501557
$ws.Range("D10").Value = 2.736028435600074
$ws.Range("E10").Value = 5.236159723822857
$ws.Range("F10").Value = 50.28161935181188
$ws.Range("G10").Value = 2.148079680861431
$ws.Range("H10").Value = 3.894742576755437
$ws.Range("I10").Value = 4.691150357187998
$ws.Range("K10").Value = 34.71931693963445
$ws.Range("L10").Value = 5.933421385815
$ws.Range("M10").Value = 13.82735339307424
$ws.Range("N10").Value = 7.046905470679085
$ws.Range("B11").Value = 16.38635543494527
$ws.Range("C11").Value = 9.676590718450644
$ws.Range("D11").Value = 2.588019284301515
$ws.Range("E11").Value = 4.892605576882899
$ws.Range("F11").Value = 47.09154245500039
$ws.Range("G11").Value = 2.145830647932072
$ws.Range("H11").Value = 4.420771876451876
$ws.Range("I11").Value = 4.622439227624828
$ws.Range("K11").Value = 32.38144840936152
$ws.Range("L11").Value = 5.656797996917494
$ws.Range("M11").Value = 14.26223436433685
$ws.Range("N11").Value = 6.588198027854022
$ws.Range("B12").Value = 16.58581023987928
$ws.Range("C12").Value = 9.792599438368384
$ws.Range("D12").Value = 2.502879083692704
$ws.Range("E12").Value = 4.652302086743603
$ws.Range("F12").Value = 44.12263087860888
$ws.Range("G12").Value = 2.14568995630251
$ws.Range("H12").Value = 5.414797094529082
$ws.Range("I12").Value = 4.623512552394335
$ws.Range("K12").Value = 30.30462102123062
$ws.Range("L12").Value = 5.473472103992671
$ws.Range("M12").Value = 14.42541469978593
$ws.Range("N12").Value = 6.149165578410881
$ws.Range("B13").Value = 16.55280963159793
$ws.Range("C13").Value = 9.747369612454101
$ws.Range("D13").Value = 2.454389208616109
$ws.Range("E13").Value = 4.486180458054197
$ws.Range("F13").Value = 40.95727276003243
$ws.Range("G13").Value = 2.14723053230089
$ws.Range("H13").Value = 6.612678149569752
$ws.Range("I13").Value = 4.690919042812937
$ws.Range("K13").Value = 28.17414931630886
$ws.Range("L13").Value = 5.348538981817542
$ws.Range("M13").Value = 14.39920538214458
$ws.Range("N13").Value = 5.704742376659531
$ws.Range("B14").Value = 16.42376470649363
$ws.Range("C14").Value = 9.646729203769372
$ws.Range("D14").Value = 2.439934951051
$ws.Range("E14").Value = 4.423232164717355
$ws.Range("F14").Value = 38.6204379872861
$ws.Range("G14").Value = 2.149065682725448
$ws.Range("H14").Value = 7.519403798601782
$ws.Range("I14").Value = 4.769866402530083
$ws.Range("K14").Value = 26.64707876589592
$ws.Range("L14").Value = 5.29978856791401
$ws.Range("M14").Value = 14.29620629081158
$ws.Range("N14").Value = 5.396899608437651
$ws.Range("B15").Value = 16.34156733274705
$ws.Range("C15").Value = 9.593026343791106
$ws.Range("D15").Value = 2.441623015274382
$ws.Range("E15").Value = 4.409797318354996
$ws.Range("F15").Value = 37.97485523867606
$ws.Range("G15").Value = 2.14990684055509
$ws.Range("H15").Value = 7.741246734921904
$ws.Range("I15").Value = 4.807562457698273
$ws.Range("K15").Value = 26.2372553781788
$ws.Range("L15").Value = 5.289940996203997
$ws.Range("M15").Value = 14.23209173926189
$ws.Range("N15").Value = 5.321007033379163
$ws.Range("B16").Value = 15.86177091442983
$ws.Range("C16").Value = 9.265578286723333
$ws.Range("D16").Value = 2.468305247646937
$ws.Range("E16").Value = 4.33936816286698
$ws.Range("F16").Value = 37.59473978960753
$ws.Range("G16").Value = 2.153653133477035
$ws.Range("H16").Value = 7.635817553456825
$ws.Range("I16").Value = 4.960928992341297
$ws.Range("K16").Value = 26.09095728237554
$ws.Range("L16").Value = 5.250371776410458
$ws.Range("M16").Value = 13.84695967969189
$ws.Range("N16").Value = 5.290636157635515
$ws.Range("B17").Value = 15.55751725581734
$ws.Range("C17").Value = 9.066097204703141
$ws.Range("D17").Value = 2.494853172694248
$ws.Range("E17").Value = 4.307417976404516
$ws.Range("F17").Value = 38.56266525047249
$ws.Range("G17").Value = 2.155567588138582
$ws.Range("H17").Value = 7.058902094856514
$ws.Range("I17").Value = 5.037685799035474
$ws.Range("K17").Value = 26.79038341206803
$ws.Range("L17").Value = 5.237981595582334
$ws.Range("M17").Value = 13.60329724743253
$ws.Range("N17").Value = 5.430383016454793
$ws.Range("B18").Value = 15.37659356468405
$ws.Range("C18").Value = 8.944847435976511
$ws.Range("D18").Value = 2.536513478781805
$ws.Range("E18").Value = 4.368576446488375
$ws.Range("F18").Value = 40.90562141893971
$ws.Range("G18").Value = 2.155998619689209
$ws.Range("H18").Value = 6.078375175361223
$ws.Range("I18").Value = 5.04550582585331
$ws.Range("K18").Value = 28.39148276123505
$ws.Range("L18").Value = 5.292212357210464
$ws.Range("M18").Value = 13.45303434845531
$ws.Range("N18").Value = 5.743332568645741
$ws.Range("B19").Value = 15.30437975571522
$ws.Range("C19").Value = 8.924065142825473
$ws.Range("D19").Value = 2.604682584188303
$ws.Range("E19").Value = 4.587609007125078
$ws.Range("F19").Value = 44.04355099618621
$ws.Range("G19").Value = 2.155150094174224
$ws.Range("H19").Value = 5.017281277315452
$ws.Range("I19").Value = 5.007293528901142
$ws.Range("K19").Value = 30.52262308121558
$ws.Range("L19").Value = 5.45351768032011
$ws.Range("M19").Value = 13.39511593204428
$ws.Range("N19").Value = 6.197452066077483
$ws.Range("B20").Value = 15.55374355830482
$ws.Range("C20").Value = 9.138163527627254
$ws.Range("D20").Value = 2.735696300706736
$ws.Range("E20").Value = 5.178898918030015
$ws.Range("F20").Value = 49.46244057788798
$ws.Range("G20").Value = 2.150672270860464
$ws.Range("H20").Value = 4.01632153637446
$ws.Range("I20").Value = 4.815359784751444
$ws.Range("K20").Value = 34.19349709607609
$ws.Range("L20").Value = 5.88312735467849
$ws.Range("M20").Value = 13.59172932276149
$ws.Range("N20").Value = 6.990882122606283
$ws.Range("B21").Value = 16.42689229349913
$ws.Range("C21").Value = 9.752370433861765
$ws.Range("D21").Value = 2.723812664641616
$ws.Range("E21").Value = 5.430673380479622
$ws.Range("F21").Value = 51.71655009559931
$ws.Range("G21").Value = 2.143149009681699
$ws.Range("H21").Value = 3.637488225404903
$ws.Range("I21").Value = 4.508476726085867
$ws.Range("K21").Value = 35.57194056249094
$ws.Range("L21").Value = 6.055922176436958
$ws.Range("M21").Value = 14.29615484362241
$ws.Range("N21").Value = 7.291704513692905
$ws.Range("B22").Value = 16.9809187090488
$ws.Range("C22").Value = 10.1258303625024
$ws.Range("D22").Value = 2.711990179802103
$ws.Range("E22").Value = 5.56202543876692
$ws.Range("F22").Value = 53.02980917425303
$ws.Range("G22").Value = 2.138374566330777
$ws.Range("H22").Value = 3.405617110969949
$ws.Range("I22").Value = 4.306570458737373
$ws.Range("K22").Value = 36.37523816588787
$ws.Range("L22").Value = 6.14949351397093
$ws.Range("M22").Value = 14.77175660826698
$ws.Range("N22").Value = 7.444047821225179
$ws.Range("B23").Value = 16.69011725154653
$ws.Range("C23").Value = 9.912272241396476
$ws.Range("D23").Value = 2.723273712979033
$ws.Range("E23").Value = 5.49158451614924
$ws.Range("F23").Value = 52.43572640253129
$ws.Range("G23").Value = 2.140869941255448
$ws.Range("H23").Value = 3.526717406530042
$ws.Range("I23").Value = 4.404140874932673
$ws.Range("K23").Value = 36.02959875029659
$ws.Range("L23").Value = 6.104613194043764
$ws.Range("M23").Value = 14.50124314832772
$ws.Range("N23").Value = 7.357288836268619
$ws.Range("B24").Value = 15.54253117482268
$ws.Range("C24").Value = 9.104729975350766
$ws.Range("D24").Value = 2.755663925040936
$ws.Range("E24").Value = 5.22071005483925
$ws.Range("F24").Value = 49.95695909254525
$ws.Range("G24").Value = 2.150538145328075
$ws.Range("H24").Value = 3.996775779812665
$ws.Range("I24").Value = 4.796810256545523
$ws.Range("K24").Value = 34.55417139379347
$ws.Range("L24").Value = 5.922981730955878
$ws.Range("M24").Value = 13.57276061628644
$ws.Range("N24").Value = 7.031877477012776
$ws.Range("B25").Value = 14.21732720277296
$ws.Range("C25").Value = 8.183437292718461
$ws.Range("D25").Value = 2.786094742153774
$ws.Range("E25").Value = 4.917512419200087
$ws.Range("F25").Value = 47.15954570897042
$ws.Range("G25").Value = 2.161302005950794
$ws.Range("H25").Value = 4.521639108364349
$ws.Range("I25").Value = 5.236524019207363
$ws.Range("K25").Value = 32.89271303371757
$ws.Range("L25").Value = 5.718929564374157
$ws.Range("M25").Value = 12.51709535426213
$ws.Range("N25").Value = 6.668833168056334
